$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").EntireRow.Delete()
$ws.Range("B16").EntireRow.Delete()

Write-Host "After deleting row 18 and row 16:"
Write-Host "B16: $($ws.Range("B16").Value())"
Write-Host "C16: $($ws.Range("C16").Value())"
Write-Host "D16: $($ws.Range("D16").Value())"
Write-Host "E16: $($ws.Range("E16").Value())"
Write-Host "F16: $($ws.Range("F16").Value())"
Write-Host "G16: $($ws.Range("G16").Value())"
Write-Host "B17: $($ws.Range("B17").Value())"
Write-Host "B21: $($ws.Range("B21").Value())"
Write-Host "B22: $($ws.Range("B22").Value())"
Write-Host "H21: $($ws.Range("H21").Value())"
Write-Host "H22: $($ws.Range("H22").Value())"
Write-Host "Dimension: $($ws.UsedRange.Address())"
